$d = $word.ActiveDocument

$d.Content.Find.Execute("`${order_Date}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${order_date}", 2)
